$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 349372.06
$ws.Range("I51").Value = 530996.75
$ws.Range("J51").Value = 2634
$ws.Range("K51").Value = 530996.75
$ws.Range("L51").Value = 2634
$ws.Range("M51").Value = -530512.75
$ws.Range("N51").Value = -3602

$ws.Range("H82").Value = 8975.111000000001
$ws.Range("I82").Value = 400
$ws.Range("J82").Value = 10047
$ws.Range("K82").Value = 1200
$ws.Range("L82").Value = 30141
$ws.Range("M82").Value = -794
$ws.Range("N82").Value = -30953

$ws.Range("H85").Value = 8975.111000000001
$ws.Range("I85").Value = 400
$ws.Range("J85").Value = 10047
$ws.Range("K85").Value = 1200
$ws.Range("L85").Value = 30141
$ws.Range("M85").Value = 204
$ws.Range("N85").Value = -32949

$ws.Range("H132").Value = 1716563.6
$ws.Range("I132").Value = 2218165.5
$ws.Range("J132").Value = 2757.5
$ws.Range("K132").Value = 6654496.5
$ws.Range("L132").Value = 8272.5
$ws.Range("M132").Value = -6651966.5
$ws.Range("N132").Value = -13332.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1848.15
$ws.Range("I32").Value = 1321.6989
$ws.Range("J32").Value = 4418.4707
$ws.Range("K32").Value = 1321.6989
$ws.Range("L32").Value = 4418.4707
$ws.Range("M32").Value = -1034.6989
$ws.Range("N32").Value = -4992.4707

$ws.Range("H61").Value = 1128.697
$ws.Range("I61").Value = 967.8889
$ws.Range("J61").Value = 1852.3334
$ws.Range("K61").Value = 967.8889
$ws.Range("L61").Value = 1852.3334
$ws.Range("M61").Value = -755.8889
$ws.Range("N61").Value = -2276.3334

$ws.Range("H92").Value = 20150
$ws.Range("J92").Value = 20150
$ws.Range("L92").Value = 20150
$ws.Range("N92").Value = -25142

$ws.Range("H102").Value = 1482.3529
$ws.Range("I102").Value = 1482.3529
$ws.Range("K102").Value = 1482.3529
$ws.Range("M102").Value = 139.6470999999999

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws.Range("H109").Value = 20377
$ws.Range("J109").Value = 20377
$ws.Range("L109").Value = 20377
$ws.Range("N109").Value = -23151

$ws.Range("H122").Value = 1988.4681
$ws.Range("I122").Value = 1311.3077
$ws.Range("J122").Value = 2826.8572
$ws.Range("K122").Value = 3933.9231
$ws.Range("L122").Value = 8480.571599999999
$ws.Range("M122").Value = -1483.9231
$ws.Range("N122").Value = -13380.5716

$ws.Range("H136").Value = 1128.697
$ws.Range("I136").Value = 967.8889
$ws.Range("J136").Value = 1852.3334
$ws.Range("K136").Value = 2903.6667
$ws.Range("L136").Value = 5557.0002
$ws.Range("M136").Value = -353.6667000000002
$ws.Range("N136").Value = -10657.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 432.22223
$ws.Range("I22").Value = 373.75
$ws.Range("K22").Value = 373.75
$ws.Range("M22").Value = -200.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3313.842
$ws.Range("I16").Value = 1868.6666
$ws.Range("J16").Value = 5791.2856
$ws.Range("K16").Value = 1868.6666
$ws.Range("L16").Value = 5791.2856
$ws.Range("M16").Value = -1581.6666
$ws.Range("N16").Value = -6365.2856

$ws.Range("H113").Value = 3313.842
$ws.Range("I113").Value = 1868.6666
$ws.Range("J113").Value = 5791.2856
$ws.Range("K113").Value = 1868.6666
$ws.Range("L113").Value = 5791.2856
$ws.Range("M113").Value = 301.3334
$ws.Range("N113").Value = -10131.2856

$ws.Range("H122").Value = 862.34485
$ws.Range("I122").Value = 683.76
$ws.Range("J122").Value = 1978.5
$ws.Range("K122").Value = 2051.28
$ws.Range("L122").Value = 5935.5
$ws.Range("M122").Value = 398.7200000000003
$ws.Range("N122").Value = -10835.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3322650
$ws.Range("I113").Value = 365.33334
$ws.Range("J113").Value = 7519220
$ws.Range("K113").Value = 1096.00002
$ws.Range("L113").Value = 22557660
$ws.Range("M113").Value = 1073.99998
$ws.Range("N113").Value = -22562000

$ws.Range("H131").Value = 759.1163
$ws.Range("I131").Value = 359.85715
$ws.Range("J131").Value = 951.86206
$ws.Range("K131").Value = 1079.57145
$ws.Range("L131").Value = 2855.58618
$ws.Range("M131").Value = 3960.42855
$ws.Range("N131").Value = -12935.58618

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 19326
$ws.Range("J123").Value = 19326
$ws.Range("L123").Value = 19326
$ws.Range("N123").Value = -24226

$ws.Range("H134").Value = 22146.572
$ws.Range("J134").Value = 22146.572
$ws.Range("L134").Value = 66439.716
$ws.Range("N134").Value = -71509.716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2125.1304
$ws.Range("I61").Value = 1535.5625
$ws.Range("J61").Value = 3472.7144
$ws.Range("K61").Value = 1535.5625
$ws.Range("L61").Value = 3472.7144
$ws.Range("M61").Value = -1333.5625
$ws.Range("N61").Value = -3876.7144

$ws.Range("H100").Value = 1800.2307
$ws.Range("I100").Value = 1286.1428
$ws.Range("K100").Value = 1286.1428
$ws.Range("M100").Value = -745.1428000000001

$ws.Range("H113").Value = 2125.1304
$ws.Range("I113").Value = 1535.5625
$ws.Range("J113").Value = 3472.7144
$ws.Range("K113").Value = 1535.5625
$ws.Range("L113").Value = 3472.7144
$ws.Range("M113").Value = 634.4375
$ws.Range("N113").Value = -7812.7144

$ws.Range("H122").Value = 5577.3784
$ws.Range("I122").Value = 5940.5806
$ws.Range("J122").Value = 3700.8333
$ws.Range("K122").Value = 17821.7418
$ws.Range("L122").Value = 11102.4999
$ws.Range("M122").Value = -15371.7418
$ws.Range("N122").Value = -16002.4999

$ws.Range("H133").Value = 23800
$ws.Range("J133").Value = 23800
$ws.Range("L133").Value = 23800
$ws.Range("N133").Value = -28860

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 722.125
$ws.Range("I113").Value = 722.125
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2166.375
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 3.625
$ws.Range("N113").ClearContents()

$ws.Range("H136").Value = 4863.4287
$ws.Range("I136").Value = 1261.3684
$ws.Range("J136").Value = 12467.777
$ws.Range("K136").Value = 3784.1052
$ws.Range("L136").Value = 37403.331
$ws.Range("M136").Value = -1234.1052
$ws.Range("N136").Value = -42503.331

